$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number (e.g. "64.99") are
# force-written as Text so Excel does not silently convert them to floats,
# then the style is reset back to Normal so no stray number-format survives.
$ws.Range("D2").Value = '26.624.75'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.596.23'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0619'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '1.820.19'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '1.594.48'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '26.602.14'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("E21").Value = '  +4.90%  '
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("D34").Value = '1.277.72'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("E35").Value = '  -8.08%  '
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("E40").Value = '  +18.59%  '
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("D45").Value = '1.732.78'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("E48").Value = '  +3.94%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.77%  '
